$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued columns (D = Price, can look numeric) need NumberFormat "@" forced
# then cleared back to default after assignment, so the resulting cell keeps
# plain text content (matching original inlineStr/shared-string text) without
# leaving a residual cell style applied.

$ws.Range("E2").Value = '  +2.68%  '
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.361.71'
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.109.40'
$ws.Range("D3").ClearFormats()

$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.70'
$ws.Range("D5").ClearFormats()

$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").ClearFormats()

$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5249'
$ws.Range("D7").ClearFormats()

$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4443'
$ws.Range("D8").ClearFormats()

$ws.Range("E9").Value = '  +4.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.74'
$ws.Range("D9").ClearFormats()

$ws.Range("E10").Value = '  +4.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09525'
$ws.Range("D10").ClearFormats()

$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.175'
$ws.Range("D11").ClearFormats()

$ws.Range("E12").Value = '  +1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.13'
$ws.Range("D12").ClearFormats()

$ws.Range("E13").Value = '  +8.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.789'
$ws.Range("D13").ClearFormats()

$ws.Range("E14").Value = '  +3.03%  '

$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.091.69'
$ws.Range("D15").ClearFormats()

$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '101.91'
$ws.Range("D16").ClearFormats()

$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001170'
$ws.Range("D17").ClearFormats()

$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.008'
$ws.Range("D18").ClearFormats()

$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.26'
$ws.Range("D19").ClearFormats()

$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06738'
$ws.Range("D20").ClearFormats()

$ws.Range("E21").Value = '  +2.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.323'
$ws.Range("D21").ClearFormats()

$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.007'
$ws.Range("D22").ClearFormats()

$ws.Range("E23").Value = '  +2.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.402.54'
$ws.Range("D23").ClearFormats()

$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.63'
$ws.Range("D24").ClearFormats()

$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.317'
$ws.Range("D25").ClearFormats()

$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.362.08'
$ws.Range("D26").ClearFormats()

$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.08'
$ws.Range("D27").ClearFormats()

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.94'
$ws.Range("D28").ClearFormats()

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.547'
$ws.Range("D29").ClearFormats()

$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.56'
$ws.Range("D30").ClearFormats()

$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.149'
$ws.Range("D31").ClearFormats()

$ws.Range("E32").Value = '  +7.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.745'
$ws.Range("D32").ClearFormats()

$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1057'
$ws.Range("D33").ClearFormats()

$ws.Range("E34").Value = '  +15.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.930'
$ws.Range("D34").ClearFormats()

$ws.Range("E35").Value = '  +2.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.270'
$ws.Range("D35").ClearFormats()

$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.932'
$ws.Range("D36").ClearFormats()

$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.53'
$ws.Range("D37").ClearFormats()

$ws.Range("E38").Value = '  +3.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02636'
$ws.Range("D38").ClearFormats()

$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06819'
$ws.Range("D39").ClearFormats()

$ws.Range("E40").Value = '  +3.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7044'
$ws.Range("D40").ClearFormats()

$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.60'
$ws.Range("D41").ClearFormats()

$ws.Range("E42").Value = '  +5.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.348'
$ws.Range("D42").ClearFormats()

$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2235'
$ws.Range("D43").ClearFormats()

$ws.Range("E44").Value = '  +3.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6854'
$ws.Range("D44").ClearFormats()

$ws.Range("E45").Value = '  +2.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.47'
$ws.Range("D45").ClearFormats()

$ws.Range("E46").Value = '  +3.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.363'
$ws.Range("D46").ClearFormats()

$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("D47").ClearFormats()

$ws.Range("E48").Value = '  +15.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.357'
$ws.Range("D48").ClearFormats()

$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.658'
$ws.Range("D49").ClearFormats()

$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000344'
$ws.Range("D50").ClearFormats()

$ws.Range("E51").Value = '  +0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.224'
$ws.Range("D51").ClearFormats()
